$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AprilRaw")

# Header row
$ws.Range("A1").Value = "Library"
$ws.Range("B1").Value = "Items owned by this library checked out at this library this month"
$ws.Range("C1").Value = "Items owned by other libraries checked out at this library this month"
$ws.Range("D1").Value = "Total circulation this month"

$ws.Range("A2").Value = "Atchison Public Library"
$ws.Range("B2").Value = 4067
$ws.Range("C2").Value = 1390
$ws.Range("D2").Value = 5457
$ws.Range("A3").Value = "Baldwin City Public Library"
$ws.Range("B3").Value = 2599
$ws.Range("C3").Value = 524
$ws.Range("D3").Value = 3123
$ws.Range("A4").Value = "Basehor Community Library"
$ws.Range("B4").Value = 7350
$ws.Range("C4").Value = 1129
$ws.Range("D4").Value = 8479
$ws.Range("A5").Value = "Bern Community Library"
$ws.Range("B5").Value = 88
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 118
$ws.Range("A6").Value = "Bonner Springs City Library"
$ws.Range("B6").Value = 4809
$ws.Range("C6").Value = 1132
$ws.Range("D6").Value = 5941
$ws.Range("A7").Value = "Burlingame Community Library"
$ws.Range("B7").Value = 437
$ws.Range("C7").Value = 211
$ws.Range("D7").Value = 648
$ws.Range("A8").Value = "Carbondale City Library"
$ws.Range("B8").Value = 419
$ws.Range("C8").Value = 117
$ws.Range("D8").Value = 536
$ws.Range("A9").Value = "Centralia Community Library"
$ws.Range("B9").Value = 220
$ws.Range("C9").Value = 46
$ws.Range("D9").Value = 266
$ws.Range("A10").Value = "Corning City Library"
$ws.Range("B10").Value = 23
$ws.Range("D10").Value = 23
$ws.Range("A11").Value = "Digital Content"
$ws.Range("A12").Value = "Doniphan County Library - Elwood"
$ws.Range("B12").Value = 75
$ws.Range("C12").Value = 20
$ws.Range("D12").Value = 95
$ws.Range("A13").Value = "Doniphan County Library - Highland"
$ws.Range("B13").Value = 144
$ws.Range("C13").Value = 77
$ws.Range("D13").Value = 221
$ws.Range("A14").Value = "Doniphan County Library - Troy"
$ws.Range("B14").Value = 461
$ws.Range("C14").Value = 254
$ws.Range("D14").Value = 715
$ws.Range("A15").Value = "Doniphan County Library - Wathena"
$ws.Range("B15").Value = 301
$ws.Range("C15").Value = 56
$ws.Range("D15").Value = 357
$ws.Range("A16").Value = "Effingham Community Library"
$ws.Range("B16").Value = 253
$ws.Range("C16").Value = 66
$ws.Range("D16").Value = 319
$ws.Range("A17").Value = "Eudora Community Library"
$ws.Range("B17").Value = 1475
$ws.Range("C17").Value = 620
$ws.Range("D17").Value = 2095
$ws.Range("A18").Value = "Everest, Barnes Reading Room"
$ws.Range("B18").Value = 78
$ws.Range("C18").Value = 71
$ws.Range("D18").Value = 149
$ws.Range("A19").Value = "Hiawatha, Morrill Public Library"
$ws.Range("B19").Value = 1575
$ws.Range("C19").Value = 557
$ws.Range("D19").Value = 2132
$ws.Range("A20").Value = "Highland Community College"
$ws.Range("B20").Value = 48
$ws.Range("C20").Value = 33
$ws.Range("D20").Value = 81
$ws.Range("A21").Value = "Holton, Beck-Bookman Library"
$ws.Range("B21").Value = 1680
$ws.Range("C21").Value = 507
$ws.Range("D21").Value = 2187
$ws.Range("A22").Value = "Horton Public Library"
$ws.Range("B22").Value = 198
$ws.Range("C22").Value = 90
$ws.Range("D22").Value = 288
$ws.Range("A23").Value = "Lansing Community Library"
$ws.Range("B23").Value = 2004
$ws.Range("C23").Value = 616
$ws.Range("D23").Value = 2620
$ws.Range("A24").Value = "Leavenworth Public Library"
$ws.Range("B24").Value = 8467
$ws.Range("C24").Value = 1710
$ws.Range("D24").Value = 10177
$ws.Range("A25").Value = "Linwood Community Library"
$ws.Range("B25").Value = 563
$ws.Range("C25").Value = 203
$ws.Range("D25").Value = 766
$ws.Range("A26").Value = "Louisburg Library"
$ws.Range("A27").Value = "Lyndon Carnegie Library"
$ws.Range("B27").Value = 253
$ws.Range("C27").Value = 201
$ws.Range("D27").Value = 454
$ws.Range("A28").Value = "McLouth Public Library"
$ws.Range("B28").Value = 255
$ws.Range("C28").Value = 94
$ws.Range("D28").Value = 349
$ws.Range("A29").Value = "Meriden-Ozawkie Public Library"
$ws.Range("B29").Value = 1430
$ws.Range("C29").Value = 571
$ws.Range("D29").Value = 2001
$ws.Range("A30").Value = "Northeast Kansas Library System"
$ws.Range("B30").Value = 14
$ws.Range("C30").Value = 38
$ws.Range("D30").Value = 52
$ws.Range("A31").Value = "Nortonville Public Library"
$ws.Range("B31").Value = 268
$ws.Range("C31").Value = 78
$ws.Range("D31").Value = 346
$ws.Range("A32").Value = "Osage City Library"
$ws.Range("B32").Value = 1277
$ws.Range("C32").Value = 371
$ws.Range("D32").Value = 1648
$ws.Range("A33").Value = "Osawatomie Public Library"
$ws.Range("B33").Value = 858
$ws.Range("C33").Value = 443
$ws.Range("D33").Value = 1301
$ws.Range("A34").Value = "Oskaloosa Public Library"
$ws.Range("B34").Value = 514
$ws.Range("C34").Value = 193
$ws.Range("D34").Value = 707
$ws.Range("A35").Value = "Ottawa Library"
$ws.Range("B35").Value = 5769
$ws.Range("C35").Value = 787
$ws.Range("D35").Value = 6556
$ws.Range("A36").Value = "Overbrook Public Library"
$ws.Range("B36").Value = 685
$ws.Range("C36").Value = 172
$ws.Range("D36").Value = 857
$ws.Range("A37").Value = "Paola Free Library"
$ws.Range("B37").Value = 3134
$ws.Range("C37").Value = 485
$ws.Range("D37").Value = 3619
$ws.Range("A38").Value = "Perry-Lecompton Community Library"
$ws.Range("B38").Value = 172
$ws.Range("C38").Value = 16
$ws.Range("D38").Value = 188
$ws.Range("A39").Value = "Pomona Community Library"
$ws.Range("B39").Value = 67
$ws.Range("C39").Value = 60
$ws.Range("D39").Value = 127
$ws.Range("A40").Value = "Prairie Hills Schools - Axtell Public School"
$ws.Range("B40").Value = 384
$ws.Range("C40").Value = 15
$ws.Range("D40").Value = 399
$ws.Range("A41").Value = "Prairie Hills Schools - Sabetha Elementary School"
$ws.Range("B41").Value = 2224
$ws.Range("C41").Value = 72
$ws.Range("D41").Value = 2296
$ws.Range("A42").Value = "Prairie Hills Schools - Sabetha High School"
$ws.Range("B42").Value = 31
$ws.Range("C42").Value = 10
$ws.Range("D42").Value = 41
$ws.Range("A43").Value = "Prairie Hills Schools - Sabetha Middle School"
$ws.Range("B43").Value = 142
$ws.Range("C43").Value = 8
$ws.Range("D43").Value = 150
$ws.Range("A44").Value = "Richmond Public Library"
$ws.Range("B44").Value = 343
$ws.Range("C44").Value = 75
$ws.Range("D44").Value = 418
$ws.Range("A45").Value = "Rossville Community Library"
$ws.Range("B45").Value = 1364
$ws.Range("C45").Value = 606
$ws.Range("D45").Value = 1970
$ws.Range("A46").Value = "Sabetha, Mary Cotton Library"
$ws.Range("B46").Value = 2435
$ws.Range("C46").Value = 1122
$ws.Range("D46").Value = 3557
$ws.Range("A47").Value = "Seneca Free Library"
$ws.Range("B47").Value = 1478
$ws.Range("C47").Value = 233
$ws.Range("D47").Value = 1711
$ws.Range("A48").Value = "Silver Lake Library"
$ws.Range("B48").Value = 889
$ws.Range("C48").Value = 658
$ws.Range("D48").Value = 1547
$ws.Range("A49").Value = "Tonganoxie Public Library"
$ws.Range("B49").Value = 2627
$ws.Range("C49").Value = 614
$ws.Range("D49").Value = 3241
$ws.Range("A50").Value = "Valley Falls, Delaware Township Library"
$ws.Range("B50").Value = 483
$ws.Range("C50").Value = 163
$ws.Range("D50").Value = 646
$ws.Range("A51").Value = "Wellsville City Library"
$ws.Range("B51").Value = 881
$ws.Range("C51").Value = 371
$ws.Range("D51").Value = 1252
$ws.Range("A52").Value = "Wetmore Public Library"
$ws.Range("B52").Value = 87
$ws.Range("C52").Value = 163
$ws.Range("D52").Value = 250
$ws.Range("A53").Value = "Williamsburg Community Library"
$ws.Range("B53").Value = 174
$ws.Range("C53").Value = 51
$ws.Range("D53").Value = 225
$ws.Range("A54").Value = "Winchester Public Library"
$ws.Range("B54").Value = 361
$ws.Range("C54").Value = 365
$ws.Range("D54").Value = 726
